# New weekly price record for "Pepino ensalada" at Vega Monumental Concepción.
# A new row of data is inserted at row 101 (the sheet is sorted with the
# newest date inserted at the top of the data block), pushing the existing
# rows 101-213 down to 102-214 and growing the used range to A1:R214.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 101, shifting rows 101:213
# down to 102:214 (mirrors Excel's own "Insert Sheet Rows").
$ws.Rows(101).Insert()

# Populate the newly inserted row 101 with the new record.
$ws.Range("A101").Value = 11
$ws.Range("B101").Value = "Vega Monumental Concepción"
$ws.Range("C101").Value = "Bíobío"
$ws.Range("D101").Value = 45079
$ws.Range("E101").Value = 8
$ws.Range("F101").Value = 100112043
$ws.Range("G101").Value = "Pepino ensalada"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 100
$ws.Range("K101").Value = 11000
$ws.Range("L101").Value = 12000
$ws.Range("M101").Value = 11500
$ws.Range("N101").Value = "$/caja 60 unidades"
$ws.Range("O101").Value = "Región de Arica y Parinacota"
$ws.Range("P101").Value = 192
$ws.Range("Q101").Value = 60
$ws.Range("R101").Value = "Hortaliza"
